# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) for a batch of Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Values below reflect the latest pulled market averages.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1099
$ws.Range("I43").Value = 1049
$ws.Range("K43").Value = 1049
$ws.Range("M43").Value = -980

# Row 107
$ws.Range("H107").Value = 112.35
$ws.Range("I107").Value = 117.75
$ws.Range("J107").Value = 90.75
$ws.Range("K107").Value = 117.75
$ws.Range("L107").Value = 90.75
$ws.Range("M107").Value = 1802.25
$ws.Range("N107").Value = -3930.75

# Row 132
$ws.Range("H132").Value = 1970.7241
$ws.Range("I132").Value = 2165
$ws.Range("J132").Value = 1038.2
$ws.Range("K132").Value = 6495
$ws.Range("L132").Value = 3114.6
$ws.Range("M132").Value = -3965
$ws.Range("N132").Value = -8174.6

# Row 138
$ws.Range("H138").Value = 1559.9688
$ws.Range("I138").Value = 1081.3182
$ws.Range("J138").Value = 1810.6904
$ws.Range("K138").Value = 3243.9546
$ws.Range("L138").Value = 5432.0712
$ws.Range("M138").Value = 1896.0454
$ws.Range("N138").Value = -15712.0712

# Row 140
$ws.Range("H140").Value = 81803.25
$ws.Range("J140").Value = 81803.25
$ws.Range("L140").Value = 81803.25
$ws.Range("N140").Value = -92163.25

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 12502886
$ws.Range("J45").Value = 31252274
$ws.Range("L45").Value = 31252274
$ws.Range("N45").Value = -31253028

# Row 104
$ws.Range("H104").Value = 27778.4
$ws.Range("J104").Value = 27778.4
$ws.Range("L104").Value = 27778.4
$ws.Range("N104").Value = -34766.4

# Row 107
$ws.Range("H107").Value = 49708.285
$ws.Range("J107").Value = 49708.285
$ws.Range("L107").Value = 49708.285
$ws.Range("N107").Value = -57388.285

# Row 134
$ws.Range("H134").Value = 123994
$ws.Range("J134").Value = 123994
$ws.Range("L134").Value = 123994
$ws.Range("N134").Value = -134134

# Row 135
$ws.Range("H135").Value = 57832
$ws.Range("J135").Value = 57832
$ws.Range("L135").Value = 57832
$ws.Range("N135").Value = -67972

# Row 138
$ws.Range("H138").Value = 61666
$ws.Range("J138").Value = 61666
$ws.Range("L138").Value = 61666
$ws.Range("N138").Value = -71946

# Row 139
$ws.Range("H139").Value = 78991
$ws.Range("J139").Value = 78991
$ws.Range("L139").Value = 78991
$ws.Range("N139").Value = -89271

# Row 140
$ws.Range("H140").Value = 77995.336
$ws.Range("J140").Value = 77995.336
$ws.Range("L140").Value = 77995.336
$ws.Range("N140").Value = -88355.336

# Row 141
$ws.Range("H141").Value = 95499.164
$ws.Range("J141").Value = 86748.75
$ws.Range("L141").Value = 86748.75
$ws.Range("N141").Value = -97108.75

$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 60000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 60000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 60000
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -61372

# Row 65
$ws.Range("H65").Value = 60000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 60000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 180000
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -186864

# Row 80
$ws.Range("H80").Value = 536.5
$ws.Range("J80").Value = 687.1
$ws.Range("L80").Value = 687.1
$ws.Range("N80").Value = -2683.1

# Row 83
$ws.Range("H83").Value = 536.5
$ws.Range("J83").Value = 687.1
$ws.Range("L83").Value = 3435.5
$ws.Range("N83").Value = -13419.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1965.6552
$ws.Range("I31").Value = 1390.3636
$ws.Range("K31").Value = 1390.3636
$ws.Range("M31").Value = -1095.3636

# Row 34
$ws.Range("H34").Value = 1965.6552
$ws.Range("I34").Value = 1390.3636
$ws.Range("K34").Value = 1390.3636
$ws.Range("M34").Value = -1188.3636

# Row 132
$ws.Range("H132").Value = 2477.7778
$ws.Range("I132").Value = 2675
$ws.Range("K132").Value = 8025
$ws.Range("M132").Value = -5495

# Row 134
$ws.Range("H134").Value = 1726.0834
$ws.Range("I134").Value = 1283
$ws.Range("K134").Value = 3849
$ws.Range("M134").Value = -1314

$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Range("H59").Value = 4002
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 4002
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 12006
$ws.Range("M59").Value = $null
$ws.Range("N59").Value = -13086

$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 25130.5
$ws.Range("J39").Value = 25130.5
$ws.Range("L39").Value = 25130.5
$ws.Range("N39").Value = -26194.5

# Row 110
$ws.Range("H110").Value = 73956.63
$ws.Range("J110").Value = 73956.63
$ws.Range("L110").Value = 73956.63
$ws.Range("N110").Value = -82136.63

$ws = $wb.Worksheets.Item("LTW")
# Row 98
$ws.Range("H98").Value = 24500
$ws.Range("J98").Value = 24500
$ws.Range("L98").Value = 24500
$ws.Range("N98").Value = -30490

# Row 122
$ws.Range("H122").Value = 66670570
$ws.Range("I122").Value = 111115180
$ws.Range("K122").Value = 333345540
$ws.Range("M122").Value = -333343090

# Row 132
$ws.Range("H132").Value = 2387.68
$ws.Range("I132").Value = 1689.8334
$ws.Range("K132").Value = 5069.5002
$ws.Range("M132").Value = -2539.5002

# Row 134
$ws.Range("H134").Value = 118572.336
$ws.Range("J134").Value = 118572.336
$ws.Range("L134").Value = 118572.336
$ws.Range("N134").Value = -128712.336

# Row 136
$ws.Range("H136").Value = 3018.946
$ws.Range("I136").Value = 3558.5557
$ws.Range("J136").Value = 2507.7368
$ws.Range("K136").Value = 10675.6671
$ws.Range("L136").Value = 7523.2104
$ws.Range("M136").Value = -8125.667099999999
$ws.Range("N136").Value = -12623.2104

# Row 138
$ws.Range("H138").Value = 95877.8
$ws.Range("J138").Value = 95877.8
$ws.Range("L138").Value = 95877.8
$ws.Range("N138").Value = -106157.8

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 95483.57000000001
$ws.Range("J46").Value = 95483.57000000001
$ws.Range("L46").Value = 95483.57000000001
$ws.Range("N46").Value = -95945.57000000001

# Row 113
$ws.Range("H113").Value = 916.4783
$ws.Range("I113").Value = 493
$ws.Range("K113").Value = 1479
$ws.Range("M113").Value = 691

# Row 122
$ws.Range("H122").Value = 3199
$ws.Range("I122").Value = 3122
$ws.Range("J122").Value = 3229.8
$ws.Range("K122").Value = 9366
$ws.Range("L122").Value = 9689.400000000001
$ws.Range("M122").Value = -6916
$ws.Range("N122").Value = -14589.4

# Row 126
$ws.Range("H126").Value = 1450.1111
$ws.Range("I126").Value = 1450.1111
$ws.Range("K126").Value = 4350.3333
$ws.Range("M126").Value = -1880.3333

# Row 134
$ws.Range("H134").Value = 95483.57000000001
$ws.Range("J134").Value = 95483.57000000001
$ws.Range("L134").Value = 286450.71
$ws.Range("N134").Value = -291520.71

# Row 141
$ws.Range("H141").Value = 57645.832
$ws.Range("J141").Value = 57645.832
$ws.Range("L141").Value = 57645.832
$ws.Range("N141").Value = -68005.83199999999
